$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column L: "break_on_off"
$ws.Range("L1").Value = "break_on_off"

# Fill default value (0) for all data rows 2-73
$ws.Range("L2:L73").Value = 0

# Rows that should be 1 (break_on_off = 1)
$ws.Range("L19").Value = 1
$ws.Range("L37").Value = 1
$ws.Range("L54").Value = 1

# Match the final selection shown in the saved file (whole new column selected)
[void]$ws.Range("L1:L73").Select()
